# Update Excel data - 2024-11-22 05:16:22
# Refreshes crypto market data across all three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Top 50 Cryptocurrencies"
# Columns: A=Name, B=Symbol, C=Current Price (USD), D=Market Capitalization,
#          E=24h Trading Volume, F=Price Change (24h %)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$cryptoUpdates = @(
    @{Row=2; C=98953; D=1957819789649; E=107473015554; F=1.48841},
    @{Row=3; C=3388.15; D=408442019407; E=55538713089; F=7.89304},
    @{Row=4; C=1.001; D=131041154867; E=113718779050; F=-0.01997},
    @{Row=5; C=261.77; D=124197184685; E=14991831812; F=8.79411},
    @{Row=6; C=635.35; D=92741072509; E=2472723832; F=3.75845},
    @{Row=7; C=1.38; D=79105543405; E=17824065953; F=23.06077},
    @{Row=8; C=0.396932; D=58333135134; E=9895261537; F=2.05635},
    @{Row=9; C=1; D=38340757364; E=11657194673; F=-0.0022},
    @{Row=10; C=3387.77; D=33192243452; E=145663717; F=8.07043},
    @{Row=11; C=0.883266; D=31681558093; E=3291559390; F=10.80615},
    @{Row=12; C=0.200654; D=17337756980; E=1073293125; F=1.51007},
    @{Row=13; C=36.46; D=14919274431; E=1037618086; F=6.6755},
    @{Row=14; C=0.00002504; D=14768327270; E=1616021481; F=3.28015},
    @{Row=15; C=4024.97; D=14527190907; E=168464106; F=7.91277},
    @{Row=16; C=98944; D=14457817201; E=864190786; F=2.04243},
    @{Row=17; C=5.57; D=14197176170; E=637419627; F=3.39082},
    @{Row=18; C=3.62; D=10307866671; E=2420681201; F=1.06291},
    @{Row=19; C=497.43; D=9840002010; E=2055046626; F=-1.2951},
    @{Row=20; C=3387.27; D=9654793051; E=1123693284; F=7.91282},
    @{Row=21; D=9580085936; E=1244740994; F=4.26916},
    @{Row=22; C=0.00002135; D=8982335231; E=6858230993; F=9.116300000000001},
    @{Row=23; C=6.23; D=8964658656; E=818046188; F=8.81039},
    @{Row=24; C=0.283567; D=8548224598; E=2314091514; F=17.29955},
    @{Row=25; C=8.800000000000001; D=8158757929; E=3462321; F=3.45041},
    @{Row=26; C=5.82; D=7093087429; E=1014559086; F=4.56338},
    @{Row=27; C=90.79000000000001; D=6832406054; E=1428455782; F=4.19564},
    @{Row=28; C=12.13; D=6468582664; E=868473780; F=4.26232},
    @{Row=29; C=3594.84; D=6254075727; E=104599642; F=8.82516},
    @{Row=30; C=9.42; D=5659344661; E=858641909; F=6.26761},
    @{Row=31; C=0.19867; D=5384790926; E=121634634; F=12.78974},
    @{Row=32; C=0.9986930000000001; D=5226843770; E=16451235; F=-0.21584},
    @{Row=33; C=0.136986; D=5220754802; E=901563713; F=8.60918},
    @{Row=34; C=9.710000000000001; D=4605777921; E=274183142; F=6.81327},
    @{Row=35; C=27.99; D=4191645494; E=889580767; F=5.49776},
    @{Row=36; C=0.00005236; D=3931813183; E=1699923348; F=2.74062},
    @{Row=37; C=7.4; D=3831971457; E=436575944; F=0.13281},
    @{Row=38; C=0.150958; D=3804917099; E=151965451; F=-0.91532},
    @{Row=39; C=0.472086; D=3760627445; E=489095714; F=6.98968},
    @{Row=40; C=508.95; D=3756554105; E=286529177; F=3.44413},
    @{Row=41; C=1.003; D=3691719798; E=223910580; F=0.06253},
    @{Row=42; C=24.81; D=3575039972; E=34121584; F=2.75394},
    @{Row=43; C=1; D=3444337143; E=160350992; F=-0.07893},
    @{Row=44; C=3.41; D=3406244105; E=1283931088; F=7.07185},
    @{Row=45; C=3.74; D=3371052995; E=302252483; F=2.84241},
    @{Row=46; D=3351227825; E=489828536; F=2.79961},
    @{Row=47; C=0.789357; D=3237546477; E=1671213831; F=11.87906},
    @{Row=48; C=161.16; D=2979464717; E=86762393; F=-1.29352},
    @{Row=49; C=1.96; D=2946226428; E=376100945; F=1.27825},
    @{Row=50; C=4.71; D=2831890252; E=585172344; F=7.04621},
    @{Row=51; C=0.839365; D=2818599241; E=181518286; F=13.91175}
)

foreach ($u in $cryptoUpdates) {
    $r = $u.Row
    if ($u.ContainsKey("C")) { $ws1.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws1.Cells.Item($r, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws1.Cells.Item($r, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $ws1.Cells.Item($r, 6).Value = $u.F }
}

# Row 51 coin swapped out: OKB (okb) -> Mantle (mnt)
$ws1.Cells.Item(51, 1).Value = "Mantle"
$ws1.Cells.Item(51, 2).Value = "mnt"

# ---------------------------------------------------------------------------
# Sheet 2: "Top 5 by Market Cap"
# Columns: A=Name, B=Market Capitalization
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")

$marketCapUpdates = @(
    @{Row=2; B=1957819789649},
    @{Row=3; B=408442019407},
    @{Row=4; B=131041154867},
    @{Row=5; B=124197184685},
    @{Row=6; B=92741072509}
)

foreach ($u in $marketCapUpdates) {
    $ws2.Cells.Item($u.Row, 2).Value = $u.B
}

# ---------------------------------------------------------------------------
# Sheet 3: "Summary"
# Columns: A=Metric, B=Value
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Summary")

# Leading apostrophe forces text so "$4360.69" isn't auto-converted to a number.
$ws3.Cells.Item(2, 2).Value = "'`$4360.69"
$ws3.Cells.Item(3, 2).Value = "XRP (23.06%)"
$ws3.Cells.Item(4, 2).Value = "Bitcoin Cash (-1.30%)"
